$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Answers: Introduction to quadratic equations",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Answers: Introduction to quadratic equations",
    2)

$d.Content.Find.Execute(
    "Tom Coleman",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Tom Coleman",
    2)

$d.Content.Find.Execute(
    "Answers to questions relating to the guide on introduction to quadratic equations.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Answers to questions relating to the guide on introduction to quadratic equations.",
    2)
